$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume data (GitHub Actions scheduled update).
# D-column numeric-looking values are forced to Text via a leading
# apostrophe so Excel's COM layer stores them exactly as strings
# (matching the workbook's existing inline-string/text cells) instead
# of re-interpreting them as numbers/dates.

# Row 2
$ws.Range("D2").Value = "'29.327.86"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").Value = "'1.840.60"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "'239.43"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").Value = "'0.6218"
$ws.Range("E6").Value = "  -1.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "'0.07323"
$ws.Range("E8").Value = "  -1.23%  "

# Row 9
$ws.Range("D9").Value = "'0.2876"
$ws.Range("E9").Value = "  -0.59%  "

# Row 10
$ws.Range("D10").Value = "'24.64"
$ws.Range("E10").Value = "  -0.67%  "

# Row 11
$ws.Range("D11").Value = "'0.07726"
$ws.Range("E11").Value = "  -0.31%  "

# Row 12
$ws.Range("D12").Value = "'1.827.85"
$ws.Range("E12").Value = "  -0.42%  "

# Row 13
$ws.Range("D13").Value = "'4.936"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6593"
$ws.Range("E14").Value = "  -2.85%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001021"
$ws.Range("E15").Value = "  +0.19%  "

# Row 16
$ws.Range("D16").Value = "'81.33"
$ws.Range("E16").Value = "  -0.95%  "

# Row 17
$ws.Range("D17").Value = "'6.227"
$ws.Range("E17").Value = "  -0.60%  "

# Row 18
$ws.Range("D18").Value = "'29.326.91"
$ws.Range("E18").Value = "  -0.10%  "

# Row 19
$ws.Range("D19").Value = "'236.38"
$ws.Range("E19").Value = "  +3.03%  "

# Row 20
$ws.Range("D20").Value = "'12.17"
$ws.Range("E20").Value = "  -1.10%  "

# Row 21
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("E22").Value = "  -3.29%  "

# Row 23
$ws.Range("D23").Value = "'1.006"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("D24").Value = "'157.20"
$ws.Range("E24").Value = "  -0.92%  "

# Row 25
$ws.Range("D25").Value = "'8.403"
$ws.Range("E25").Value = "  -0.88%  "

# Row 26
$ws.Range("D26").Value = "'0.1328"
$ws.Range("E26").Value = "  -1.94%  "

# Row 27
$ws.Range("D27").Value = "'17.19"
$ws.Range("E27").Value = "  -1.52%  "

# Row 28
$ws.Range("D28").Value = "'0.06869"
$ws.Range("E28").Value = "  +3.63%  "

# Row 29
$ws.Range("E29").Value = "  +0.88%  "

# Row 30
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("D31").Value = "'4.011"
$ws.Range("E31").Value = "  -1.51%  "

# Row 32
$ws.Range("D32").Value = "'3.931"
$ws.Range("E32").Value = "  -3.53%  "

# Row 33
$ws.Range("E33").Value = "  +1.10%  "

# Row 34
$ws.Range("D34").Value = "'1.741"
$ws.Range("E34").Value = "  -5.24%  "

# Row 35
$ws.Range("D35").Value = "'0.6770"
$ws.Range("E35").Value = "  -2.40%  "

# Row 36
$ws.Range("E36").Value = "  +0.28%  "

# Row 37
$ws.Range("D37").Value = "'0.01819"
$ws.Range("E37").Value = "  -2.21%  "

# Row 38
$ws.Range("D38").Value = "'2.784"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
$ws.Range("D39").Value = "'1.231.16"
$ws.Range("E39").Value = "  -1.34%  "

# Row 40
$ws.Range("D40").Value = "'6.663"
$ws.Range("E40").Value = "  -1.65%  "

# Row 41
$ws.Range("D41").Value = "'0.9424"
$ws.Range("E41").Value = "  +0.50%  "

# Row 42
$ws.Range("E42").Value = "  +0.06%  "

# Row 43
$ws.Range("D43").Value = "'1.990.92"
$ws.Range("E43").Value = "  -1.86%  "

# Row 44
$ws.Range("D44").Value = "'101.23"
$ws.Range("E44").Value = "  +0.44%  "

# Row 45
$ws.Range("D45").Value = "'65.04"
$ws.Range("E45").Value = "  -0.88%  "

# Row 46
$ws.Range("E46").Value = "  +1.44%  "

# Row 47
$ws.Range("D47").Value = "'1.683"
$ws.Range("E47").Value = "  -1.70%  "

# Row 48
$ws.Range("D48").Value = "'6.868"
$ws.Range("E48").Value = "  -2.56%  "

# Row 49
$ws.Range("D49").Value = "'8.803"
$ws.Range("E49").Value = "  -2.40%  "

# Row 50
$ws.Range("E50").Value = "  -2.27%  "

# Row 51
$ws.Range("D51").Value = "'0.3851"
$ws.Range("E51").Value = "  -1.33%  "
